$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update threshold values for alpha, beta, and ratio rows
$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 5.5
$ws.Range("B4").Value = 0.7
$ws.Range("C4").Value = 1.3

# Delete the theta_threshold_range row (row 5), shifting rows up
$ws.Rows.Item(5).Delete()

# Update selection to match final state
$ws.Range("C5").Select()

# Match the page setup recorded for this sheet (paper size/orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
